# Connections_Table.xlsx edit
# "Added Connections table to lab report"
#
# Adds a 3rd "POT" connections mini-table (headers + 3 rows) in columns I:K
# around rows 21-25, updates the H-bridge(1) label in the existing Motor2
# table to mention the POT wiring, adds a matching "Pot(2)" / "N/A" entry,
# and moves the active selection to O39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New strings get created in this order so the shared-string table ends
# --- up appended the same way the original author typed them: POT, the J11
# --- (34) pin reference, the amended H-bridge(1) label, then Pot(2).

# New "POT" mini table header label (I21)
$ws.Range("I21").Value = "POT"

# POT table column headers (row 22), same headers as the other mini tables
$ws.Range("I22").Value = "PIN #:"
$ws.Range("J22").Value = "Connection"
$ws.Range("K22").Value = "Wire Color"

# POT table body rows (23-25)
$ws.Range("I23").Value = 1
$ws.Range("J23").Value = "Expansion Board J11 (28)"
$ws.Range("K23").Value = "Red"

$ws.Range("I24").Value = 2
$ws.Range("J24").Value = "Expansion Board J11 (34)"
$ws.Range("K24").Value = "Blue"

$ws.Range("I25").Value = 3
$ws.Range("J25").Value = "Expansion Board J11 (22)"
$ws.Range("K25").Value = "Black"

# Match the Wire Color cell fills used elsewhere in the sheet for
# Red / Blue / Black (copy format only, keep the values just written).
$ws.Range("K4").Copy()
$ws.Range("K23").PasteSpecial(-4122)

$ws.Range("C10").Copy()
$ws.Range("K24").PasteSpecial(-4122)

$ws.Range("O8").Copy()
$ws.Range("K25").PasteSpecial(-4122)

# Existing Motor2 H-bridge table: B10 label now also calls out Pot(1)
$ws.Range("B10").Value = "H-bridge (1) Pot(1)"

# Row 18 (previously "N/A"/"N/A") now documents Pot(2) wiring
$ws.Range("B18").Value = "Pot(2)"
$ws.Range("C18").Value = "Blue"
$ws.Range("C10").Copy()
$ws.Range("C18").PasteSpecial(-4122)

# Move / update the active selection like the saved workbook shows
[void]$ws.Range("O39").Select()

[void]$wb.Save()
